$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (values refreshed by scheduled GitHub Actions run).

# Row 2: update D2, E2
$ws.Range("D2").Value = "64.154.49"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3: update D3, E3
$ws.Range("D3").Value = "2.762.03"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4: update D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.20"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.66"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8: update E8
$ws.Range("E8").Value = "  -2.64%  "

# Row 9: update E9
$ws.Range("E9").Value = "  -1.47%  "

# Row 10: update E10
$ws.Range("E10").Value = "  +4.04%  "

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("E11").Value = "  -13.95%  "

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  -1.06%  "

# Row 13: update D13, E13
$ws.Range("D13").Value = "3.248.12"
$ws.Range("E13").Value = "  +0.10%  "

# Row 14: update D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.14"
$ws.Range("E14").Value = "  -2.13%  "

# Row 15: update D15
$ws.Range("D15").Value = "63.754.58"

# Row 16: update D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000153"
$ws.Range("E16").Value = "  -2.66%  "

# Row 17: update D17, E17
$ws.Range("D17").Value = "2.763.46"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18: update D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.27"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19: update E19
$ws.Range("E19").Value = "  -1.43%  "

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "359.09"
$ws.Range("E20").Value = "  -1.44%  "

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  -3.09%  "

# Row 22: update B22, C22, D22, E22
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.538"
$ws.Range("E22").Value = "  -0.73%  "

# Row 23: update B23, C23, D23, E23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.69"
$ws.Range("E24").Value = "  -1.77%  "

# Row 25: update E25
$ws.Range("E25").Value = "  -0.95%  "

# Row 26: update D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.65"
$ws.Range("E26").Value = "  +0.17%  "

# Row 27: update E27
$ws.Range("E27").Value = "  +0.05%  "

# Row 28: update D28, E28
$ws.Range("D28").Value = "0.0₃0919"
$ws.Range("E28").Value = "  +0.09%  "

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  +0.93%  "

# Row 30: update E30
$ws.Range("E30").Value = "  -2.93%  "

# Row 31: update E31
$ws.Range("E31").Value = "  -1.74%  "

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.04"
$ws.Range("E32").Value = "  -2.65%  "

# Row 33: update D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.01"
$ws.Range("E33").Value = "  +1.10%  "

# Row 34: update D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.27"
$ws.Range("E34").Value = "  -1.97%  "

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  +1.99%  "

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.06%  "

# Row 37: update E37
$ws.Range("E37").Value = "  -0.33%  "

# Row 38: update E38
$ws.Range("E38").Value = "  -1.05%  "

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.48"
$ws.Range("E39").Value = "  +4.23%  "

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "349.83"
$ws.Range("E40").Value = "  +2.91%  "

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  -1.04%  "

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.19"
$ws.Range("E42").Value = "  -1.40%  "

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.77"
$ws.Range("E43").Value = "  -0.59%  "

# Row 44: update D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.03"
$ws.Range("E44").Value = "  -2.51%  "

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0594"
$ws.Range("E45").Value = "  -2.00%  "

# Row 46: update B46, C46, D46, E46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.54"
$ws.Range("E46").Value = "  -0.44%  "

# Row 47: update B47, C47, D47, E47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0257"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48: update E48
$ws.Range("E48").Value = "  -1.91%  "

# Row 49: update E49
$ws.Range("E49").Value = "  -0.38%  "

# Row 50: update E50
$ws.Range("E50").Value = "  -0.18%  "

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  +0.04%  "

